$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows at 58..67, each inheriting the number format/style
# of the row above it (xlFormatFromLeftOrAbove) so the new cells pick up
# the same style index already used by row 57 (right-aligned Arial).
# Insert() shifts row 101+ down by one each time; that is corrected
# below by deleting the equivalent number of (now-blank) rows that were
# pushed in front of the original row 101 block.
for ($i = 0; $i -lt 10; $i++) {
    $newRow = 58 + $i
    $ws.Rows.Item($newRow).Insert(-4121, 0)
    # Drop the inherited column V cell - the new rows only use A:U
    $ws.Cells.Item($newRow, 22).Clear()
}

# Undo the cumulative 10-row downward shift applied to the untouched
# row 101+ block so it lands back on its original row numbers.
for ($i = 0; $i -lt 10; $i++) {
    $ws.Rows.Item(101).Delete()
}

# Populate the new data rows with their values
$row = 58
$ws.Cells.Item($row, 1).Value() = 0.09724102198
$ws.Cells.Item($row, 2).Value() = 0.1101479291
$ws.Cells.Item($row, 3).Value() = 0.164543884
$ws.Cells.Item($row, 4).Value() = 0.3746661022
$ws.Cells.Item($row, 5).Value() = 0.6594350099
$ws.Cells.Item($row, 6).Value() = 0.2859206902
$ws.Cells.Item($row, 7).Value() = 0.5429347518
$ws.Cells.Item($row, 8).Value() = 0.2337670974
$ws.Cells.Item($row, 9).Value() = 0.1739029009
$ws.Cells.Item($row, 10).Value() = 0.113630948
$ws.Cells.Item($row, 11).Value() = 0.6657297243
$ws.Cells.Item($row, 12).Value() = 1.0
$ws.Cells.Item($row, 13).Value() = 0.6068154292
$ws.Cells.Item($row, 14).Value() = 0.1778153883
$ws.Cells.Item($row, 15).Value() = 0.08398864763
$ws.Cells.Item($row, 16).Value() = 0.1198268512
$ws.Cells.Item($row, 17).Value() = 0.1085874461
$ws.Cells.Item($row, 18).Value() = 0.2074315516
$ws.Cells.Item($row, 19).Value() = 0.7335296893
$ws.Cells.Item($row, 20).Value() = 0.849334208
$ws.Cells.Item($row, 21).Value() = 0.2410756984

$row = 59
$ws.Cells.Item($row, 1).Value() = 0.1337662451
$ws.Cells.Item($row, 2).Value() = 0.05141706438
$ws.Cells.Item($row, 3).Value() = 0.08693464737
$ws.Cells.Item($row, 4).Value() = 0.2348312391
$ws.Cells.Item($row, 5).Value() = 0.2152395099
$ws.Cells.Item($row, 6).Value() = 0.1594663913
$ws.Cells.Item($row, 7).Value() = 0.2648267579
$ws.Cells.Item($row, 8).Value() = 0.07018314442
$ws.Cells.Item($row, 9).Value() = 0.1033953884
$ws.Cells.Item($row, 10).Value() = 0.04313906364
$ws.Cells.Item($row, 11).Value() = 0.8646910972
$ws.Cells.Item($row, 12).Value() = 0.03413221131
$ws.Cells.Item($row, 13).Value() = 0.06169695919
$ws.Cells.Item($row, 14).Value() = 0.037153804
$ws.Cells.Item($row, 15).Value() = 0.06685688028
$ws.Cells.Item($row, 16).Value() = 0.3795986651
$ws.Cells.Item($row, 17).Value() = 0.07148721113
$ws.Cells.Item($row, 18).Value() = 0.10869324
$ws.Cells.Item($row, 19).Value() = 0.1123685132
$ws.Cells.Item($row, 20).Value() = 1.0
$ws.Cells.Item($row, 21).Value() = 0.2897721959

$row = 60
$ws.Cells.Item($row, 1).Value() = 0.03417009981
$ws.Cells.Item($row, 2).Value() = 0.07989124366
$ws.Cells.Item($row, 3).Value() = 0.09648353978
$ws.Cells.Item($row, 4).Value() = 0.1234646188
$ws.Cells.Item($row, 5).Value() = 0.3462660659
$ws.Cells.Item($row, 6).Value() = 1.0
$ws.Cells.Item($row, 7).Value() = 0.770740854
$ws.Cells.Item($row, 8).Value() = 0.1441139853
$ws.Cells.Item($row, 9).Value() = 0.1031161374
$ws.Cells.Item($row, 10).Value() = 0.07296654047
$ws.Cells.Item($row, 11).Value() = 0.3067376429
$ws.Cells.Item($row, 12).Value() = 0.05815123068
$ws.Cells.Item($row, 13).Value() = 0.2041629881
$ws.Cells.Item($row, 14).Value() = 0.05848334431
$ws.Cells.Item($row, 15).Value() = 0.04489834095
$ws.Cells.Item($row, 16).Value() = 0.05916289651
$ws.Cells.Item($row, 17).Value() = 0.06708557216
$ws.Cells.Item($row, 18).Value() = 0.1911687383
$ws.Cells.Item($row, 19).Value() = 0.1463332686
$ws.Cells.Item($row, 20).Value() = 0.09515224277
$ws.Cells.Item($row, 21).Value() = 0.04980511518

$row = 61
$ws.Cells.Item($row, 1).Value() = 0.1693143913
$ws.Cells.Item($row, 2).Value() = 0.3147371522
$ws.Cells.Item($row, 3).Value() = 0.2170821589
$ws.Cells.Item($row, 4).Value() = 0.5692811302
$ws.Cells.Item($row, 5).Value() = 0.6082392716
$ws.Cells.Item($row, 6).Value() = 0.4653296623
$ws.Cells.Item($row, 7).Value() = 0.8903607773
$ws.Cells.Item($row, 8).Value() = 0.3380258688
$ws.Cells.Item($row, 9).Value() = 0.1647020397
$ws.Cells.Item($row, 10).Value() = 0.2440425838
$ws.Cells.Item($row, 11).Value() = 0.5905884821
$ws.Cells.Item($row, 12).Value() = 0.3644683064
$ws.Cells.Item($row, 13).Value() = 0.1783314054
$ws.Cells.Item($row, 14).Value() = 0.07378557335
$ws.Cells.Item($row, 15).Value() = 0.1955210717
$ws.Cells.Item($row, 16).Value() = 0.1120029162
$ws.Cells.Item($row, 17).Value() = 0.1119994046
$ws.Cells.Item($row, 18).Value() = 0.5733307319
$ws.Cells.Item($row, 19).Value() = 0.1516450319
$ws.Cells.Item($row, 20).Value() = 1.0
$ws.Cells.Item($row, 21).Value() = 0.610719566

$row = 62
$ws.Cells.Item($row, 1).Value() = 0.03022771037
$ws.Cells.Item($row, 2).Value() = 0.08402835097
$ws.Cells.Item($row, 3).Value() = 0.1419964994
$ws.Cells.Item($row, 4).Value() = 0.3193490056
$ws.Cells.Item($row, 5).Value() = 1.0
$ws.Cells.Item($row, 6).Value() = 0.1403483953
$ws.Cells.Item($row, 7).Value() = 0.4288770431
$ws.Cells.Item($row, 8).Value() = 0.1063364752
$ws.Cells.Item($row, 9).Value() = 0.1489295003
$ws.Cells.Item($row, 10).Value() = 0.213108915
$ws.Cells.Item($row, 11).Value() = 0.4593016725
$ws.Cells.Item($row, 12).Value() = 0.2313420695
$ws.Cells.Item($row, 13).Value() = 0.1348333575
$ws.Cells.Item($row, 14).Value() = 0.02456818412
$ws.Cells.Item($row, 15).Value() = 0.1044198906
$ws.Cells.Item($row, 16).Value() = 0.06922586404
$ws.Cells.Item($row, 17).Value() = 0.1683641566
$ws.Cells.Item($row, 18).Value() = 0.3515515973
$ws.Cells.Item($row, 19).Value() = 0.1895871849
$ws.Cells.Item($row, 20).Value() = 0.4905839119
$ws.Cells.Item($row, 21).Value() = 0.5175832181

$row = 63
$ws.Cells.Item($row, 1).Value() = 0.1527425764
$ws.Cells.Item($row, 2).Value() = 0.2782209816
$ws.Cells.Item($row, 3).Value() = 0.1925036367
$ws.Cells.Item($row, 4).Value() = 1.0
$ws.Cells.Item($row, 5).Value() = 0.8115155891
$ws.Cells.Item($row, 6).Value() = 0.146480129
$ws.Cells.Item($row, 7).Value() = 0.5774551316
$ws.Cells.Item($row, 8).Value() = 0.1716497233
$ws.Cells.Item($row, 9).Value() = 0.1656732263
$ws.Cells.Item($row, 10).Value() = 0.03718954806
$ws.Cells.Item($row, 11).Value() = 0.5056576536
$ws.Cells.Item($row, 12).Value() = 0.1358458501
$ws.Cells.Item($row, 13).Value() = 0.09643846694
$ws.Cells.Item($row, 14).Value() = 0.07508771035
$ws.Cells.Item($row, 15).Value() = 0.07618347503
$ws.Cells.Item($row, 16).Value() = 0.04062615141
$ws.Cells.Item($row, 17).Value() = 0.06351839377
$ws.Cells.Item($row, 18).Value() = 0.3160014368
$ws.Cells.Item($row, 19).Value() = 0.09628002634
$ws.Cells.Item($row, 20).Value() = 0.5270604612
$ws.Cells.Item($row, 21).Value() = 0.2356027875

$row = 64
$ws.Cells.Item($row, 1).Value() = 0.1896328065
$ws.Cells.Item($row, 2).Value() = 0.4427254655
$ws.Cells.Item($row, 3).Value() = 0.2554946796
$ws.Cells.Item($row, 4).Value() = 0.4861495819
$ws.Cells.Item($row, 5).Value() = 0.802071037
$ws.Cells.Item($row, 6).Value() = 0.286147183
$ws.Cells.Item($row, 7).Value() = 0.4558839223
$ws.Cells.Item($row, 8).Value() = 0.1962051411
$ws.Cells.Item($row, 9).Value() = 0.2891635116
$ws.Cells.Item($row, 10).Value() = 0.1362896737
$ws.Cells.Item($row, 11).Value() = 0.669839359
$ws.Cells.Item($row, 12).Value() = 0.4397623507
$ws.Cells.Item($row, 13).Value() = 0.1378066708
$ws.Cells.Item($row, 14).Value() = 0.0489787577
$ws.Cells.Item($row, 15).Value() = 0.168837894
$ws.Cells.Item($row, 16).Value() = 0.5627751943
$ws.Cells.Item($row, 17).Value() = 0.4372578287
$ws.Cells.Item($row, 18).Value() = 0.2410277248
$ws.Cells.Item($row, 19).Value() = 0.4346455646
$ws.Cells.Item($row, 20).Value() = 0.9375894327
$ws.Cells.Item($row, 21).Value() = 0.4170547121

$row = 65
$ws.Cells.Item($row, 1).Value() = 0.1160246202
$ws.Cells.Item($row, 2).Value() = 0.2378009565
$ws.Cells.Item($row, 3).Value() = 0.3801532064
$ws.Cells.Item($row, 4).Value() = 0.7790974334
$ws.Cells.Item($row, 5).Value() = 0.7424492624
$ws.Cells.Item($row, 6).Value() = 0.2088041316
$ws.Cells.Item($row, 7).Value() = 0.2142618405
$ws.Cells.Item($row, 8).Value() = 0.113460463
$ws.Cells.Item($row, 9).Value() = 0.3262694311
$ws.Cells.Item($row, 10).Value() = 0.1522816943
$ws.Cells.Item($row, 11).Value() = 0.3119665407
$ws.Cells.Item($row, 12).Value() = 0.1849026969
$ws.Cells.Item($row, 13).Value() = 0.5127918004
$ws.Cells.Item($row, 14).Value() = 0.06365933792
$ws.Cells.Item($row, 15).Value() = 0.1816753008
$ws.Cells.Item($row, 16).Value() = 0.1401719107
$ws.Cells.Item($row, 17).Value() = 0.1116843593
$ws.Cells.Item($row, 18).Value() = 0.2342639624
$ws.Cells.Item($row, 19).Value() = 0.3085088784
$ws.Cells.Item($row, 20).Value() = 1.0
$ws.Cells.Item($row, 21).Value() = 0.4045860426

$row = 66
$ws.Cells.Item($row, 1).Value() = 0.1477363119
$ws.Cells.Item($row, 2).Value() = 0.2872983831
$ws.Cells.Item($row, 3).Value() = 0.3102443806
$ws.Cells.Item($row, 4).Value() = 0.8002844979
$ws.Cells.Item($row, 5).Value() = 0.7122537975
$ws.Cells.Item($row, 6).Value() = 0.2237407246
$ws.Cells.Item($row, 7).Value() = 1.0
$ws.Cells.Item($row, 8).Value() = 0.2776468362
$ws.Cells.Item($row, 9).Value() = 0.1621913711
$ws.Cells.Item($row, 10).Value() = 0.1122796241
$ws.Cells.Item($row, 11).Value() = 0.5239527562
$ws.Cells.Item($row, 12).Value() = 0.2732414466
$ws.Cells.Item($row, 13).Value() = 0.3002745799
$ws.Cells.Item($row, 14).Value() = 0.1045621063
$ws.Cells.Item($row, 15).Value() = 0.1417995479
$ws.Cells.Item($row, 16).Value() = 0.1272257293
$ws.Cells.Item($row, 17).Value() = 0.3578023447
$ws.Cells.Item($row, 18).Value() = 0.2538563992
$ws.Cells.Item($row, 19).Value() = 0.2276945102
$ws.Cells.Item($row, 20).Value() = 0.3483116576
$ws.Cells.Item($row, 21).Value() = 0.1868175777

$row = 67
$ws.Cells.Item($row, 1).Value() = 0.05611562564
$ws.Cells.Item($row, 2).Value() = 0.02975984176
$ws.Cells.Item($row, 3).Value() = 0.1234927588
$ws.Cells.Item($row, 4).Value() = 0.1246238893
$ws.Cells.Item($row, 5).Value() = 1.0
$ws.Cells.Item($row, 6).Value() = 0.3479713951
$ws.Cells.Item($row, 7).Value() = 0.4641026943
$ws.Cells.Item($row, 8).Value() = 0.3718869039
$ws.Cells.Item($row, 9).Value() = 0.1460528841
$ws.Cells.Item($row, 10).Value() = 0.09859125619
$ws.Cells.Item($row, 11).Value() = 0.14032932
$ws.Cells.Item($row, 12).Value() = 0.07895808724
$ws.Cells.Item($row, 13).Value() = 0.08910917938
$ws.Cells.Item($row, 14).Value() = 0.02567254796
$ws.Cells.Item($row, 15).Value() = 0.02317807762
$ws.Cells.Item($row, 16).Value() = 0.1306214366
$ws.Cells.Item($row, 17).Value() = 0.01995593165
$ws.Cells.Item($row, 18).Value() = 0.1073468266
$ws.Cells.Item($row, 19).Value() = 0.04286177531
$ws.Cells.Item($row, 20).Value() = 0.4197412262
$ws.Cells.Item($row, 21).Value() = 0.2615420383

Write-Output "Rows 58-67 added"